# Switzerland Super League workbook update
# The source feed re-sent several fixtures with the Home/Away rows in the
# opposite order from the previous pull (same match id pairs, rows 2-at-a-
# time, occasionally a 3-row rotation), plus a handful of odds refreshes on
# still-unplayed fixtures near the bottom of the sheet.
#
# Strategy: for every affected block, snapshot the B:AC payload (id through
# PL_AhUnder) of each row involved with .Value(), then write the snapshots
# back into the rows in their new order. Column A (the running match index)
# is left untouched throughout, exactly as in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rangeA = $ws.Range("B$r1`:AC$r1")
    $rangeB = $ws.Range("B$r2`:AC$r2")
    $valA = $rangeA.Value()
    $valB = $rangeB.Value()
    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

function Rotate-Rows3($r1, $r2, $r3) {
    # new(r1) = old(r3); new(r2) = old(r1); new(r3) = old(r2)
    $rangeA = $ws.Range("B$r1`:AC$r1")
    $rangeB = $ws.Range("B$r2`:AC$r2")
    $rangeC = $ws.Range("B$r3`:AC$r3")
    $valA = $rangeA.Value()
    $valB = $rangeB.Value()
    $valC = $rangeC.Value()
    $rangeA.Value = $valC
    $rangeB.Value = $valA
    $rangeC.Value = $valB
}

# Simple two-row swaps (same match-id pair, order flipped)
Swap-Rows 50 51
Swap-Rows 72 73
Swap-Rows 97 98
Swap-Rows 99 100
Swap-Rows 124 125
Swap-Rows 130 131
Swap-Rows 154 155
Swap-Rows 164 165
Swap-Rows 170 171
Swap-Rows 184 185
Swap-Rows 195 196
Swap-Rows 197 198
Swap-Rows 214 215

# Three-row rotations
Rotate-Rows3 93 94 95
Rotate-Rows3 224 225 226

# Odds refresh on a few still-unplayed fixtures (no row reordering here)
$ws.Range("U240").Value = 1.825
$ws.Range("V240").Value = 2.025

$ws.Range("N241").Value = 1.8
$ws.Range("O241").Value = 3.8
$ws.Range("P241").Value = 4.2
$ws.Range("Q241").Value = -0.5
$ws.Range("U241").Value = 2.025
$ws.Range("V241").Value = 1.825

$ws.Range("N242").Value = 2.9
$ws.Range("O242").Value = 3.1
$ws.Range("R242").Value = 2.05
$ws.Range("S242").Value = 1.8
$ws.Range("U242").Value = 1.95
$ws.Range("V242").Value = 1.9
